$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look numeric to Excel auto-detection;
# force Text number format first so COM does not coerce them to numbers,
# then restore the default "Normal" style so no stray style index remains.
$textForceCells = @(
    "D5",
    "D7",
    "D9",
    "D10",
    "D12",
    "D14",
    "D15",
    "D16",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D35",
    "D39",
    "D40",
    "D42",
    "D44",
    "D45",
    "D48",
    "D49",
    "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all cell value updates described by the diff
$ws.Range("D2").Value = "37.279.82"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "2.028.16"
$ws.Range("E3").Value = "  +3.91%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "247.90"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("D7").Value = "60.44"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.394"
$ws.Range("E9").Value = "  +4.77%  "
$ws.Range("D10").Value = "0.0809"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("D12").Value = "15.25"
$ws.Range("E12").Value = "  +7.07%  "
$ws.Range("D13").Value = "2.331.54"
$ws.Range("E13").Value = "  +4.30%  "
$ws.Range("D14").Value = "0.857"
$ws.Range("E14").Value = "  +3.86%  "
$ws.Range("D15").Value = "22.19"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D16").Value = "5.51"
$ws.Range("E16").Value = "  +5.36%  "
$ws.Range("D17").Value = "2.032.53"
$ws.Range("E17").Value = "  +3.82%  "
$ws.Range("D18").Value = "37.213.43"
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("D19").Value = "70.62"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").Value = "0.0₃0864"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").Value = "5.24"
$ws.Range("E21").Value = "  +3.50%  "
$ws.Range("D22").Value = "230.98"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "2.61"
$ws.Range("E24").Value = "  +7.09%  "
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "9.46"
$ws.Range("E26").Value = "  +3.44%  "
$ws.Range("D27").Value = "163.21"
$ws.Range("E27").Value = "  +2.12%  "
$ws.Range("D28").Value = "0.139"
$ws.Range("E28").Value = "  -3.67%  "
$ws.Range("D29").Value = "19.78"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").Value = "1.37"
$ws.Range("E30").Value = "  +4.48%  "
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("D32").Value = "4.82"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").Value = "0.0675"
$ws.Range("E33").Value = "  +10.32%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "2.49"
$ws.Range("E34").Value = "  +9.86%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "4.48"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("E36").Value = "  +4.76%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("D39").Value = "5.45"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "3.03"
$ws.Range("E40").Value = "  +3.89%  "
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").Value = "17.11"
$ws.Range("E42").Value = "  +9.14%  "
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").Value = "0.0214"
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").Value = "92.39"
$ws.Range("E45").Value = "  +4.51%  "
$ws.Range("E46").Value = "  +3.91%  "
$ws.Range("D47").Value = "1.382.00"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").Value = "7.46"
$ws.Range("E48").Value = "  +5.04%  "
$ws.Range("D49").Value = "2.16"
$ws.Range("E49").Value = "  +18.45%  "
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").Value = "46.27"
$ws.Range("E51").Value = "  +2.07%  "

# Restore default styling on the cells we temporarily formatted as Text
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
